$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge "THU MAY 31" + " 11:43:20 IST 2018" (two runs) into a single
#    run "THU MAY 31 11:43:20 IST 2018" via Find & Replace, which
#    naturally collapses the two adjoining runs into one run.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("THU MAY 31 11:43:20 IST 2018", $false, $false, $false, $false, $false, `
              $true, 1, $false, "THU MAY 31 11:43:20 IST 2018", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Locate the final "Amount Received mode ... - CASH" paragraph (the
#    very last chick-in record in the document, which is immediately
#    followed only by trailing empty paragraphs) and append a brand new
#    chick-in record after it:
#
#      SAT Jun 02 10:57:34 IST 2018
#      Person Name                    - JAYAKKA
#      ---------------------------------------------------------------
#      Item Name                      - CARROT
#      Amount Received                - 1246        (red)
#      Amount Received mode           - CASH AND CLEARD
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Amount Received mode`t`t- CASH`r") {
        $target = $p
    }
}

$r = $target.Range.Duplicate
$r.Collapse(0)

# Blank separator paragraph before the new record.
$r.InsertAfter("`r")
$r.Collapse(0)

# "SAT Jun 02" / " 10:57:34 IST 2018" -> two runs within one paragraph.
$r.InsertAfter("SAT Jun 02")
$r.Collapse(0)
$r.InsertAfter(" 10:57:34 IST 2018")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# "Person Name" <tab><tab><tab><tab>"- JAYAKKA"
$r.InsertAfter("Person Name")
$r.Collapse(0)
$r.InsertAfter("`t")
$r.Collapse(0)
$r.InsertAfter("`t")
$r.Collapse(0)
$r.InsertAfter("`t")
$r.Collapse(0)
$r.InsertAfter("`t- JAYAKKA")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# Divider line.
$r.InsertAfter("---------------------------------------------------------------")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# "Item Name" <tab><tab><tab><tab>"- CARROT"
$r.InsertAfter("Item Name")
$r.Collapse(0)
$r.InsertAfter("`t")
$r.Collapse(0)
$r.InsertAfter("`t")
$r.Collapse(0)
$r.InsertAfter("`t")
$r.Collapse(0)
$r.InsertAfter("`t- CARROT")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# "Amount Received" <tab><tab><tab>"- 1246" in red.
$r.InsertAfter("Amount Received")
$r.Collapse(0)
$r.InsertAfter("`t")
$r.Collapse(0)
$r.InsertAfter("`t")
$r.Collapse(0)
$r.InsertAfter("`t- 1246")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# "Amount Received mode" <tab><tab>"- CASH AND CLEARD"
$r.InsertAfter("Amount Received mode")
$r.Collapse(0)
$r.InsertAfter("`t")
$r.Collapse(0)
$r.InsertAfter("`t- CASH AND CLEARD")
$r.Collapse(0)

# Trailing blank paragraph closing the new record.
$r.InsertAfter("`r")
$r.Collapse(0)

# ---------------------------------------------------------------------
# 3) Colour the "Amount Received" paragraph of the new record red
#    (FF0000), matching the other records' convention.
# ---------------------------------------------------------------------
$amountPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Amount Received`t`t`t- 1246`r") {
        $amountPara = $p
    }
}
if ($amountPara -ne $null) {
    $amountPara.Range.Font.Color = 255
}

Write-Output "done"
